$wb = $excel.ActiveWorkbook
$wsFunc = $wb.Worksheets.Item("Functional Requirements")
$wsDomain = $wb.Worksheets.Item("Domain Requirements")

# Apply text edits in the same order the original author made them, so that
# newly created shared-string entries are appended in the matching order.

# 1-2: "transaction/error logs" -> "transaction/error log" (titles)
$wsFunc.Range("C15").Value = "The TVM must be able to show transaction log"
$wsFunc.Range("C16").Value = "The TVM must be able to show error log"

# 3-4: matching description text
$wsFunc.Range("E15").Value = "The TVM operator must be able to view transaction log to see the detail of each transaction."
$wsFunc.Range("E16").Value = "The TVM operator must be able to view error log to see the detail of each error."

# 5: fix typo "datas" -> "data"
$wsDomain.Range("E9").Value = "The TVM must ensure the security of all transactions and other sensitive data"

# 6-7: "credit/debit card" -> "ATM card" on the Functional Requirements sheet
$wsFunc.Range("E6").Value = "The TVM must be able to accept many payment methods, including ATM card and mobile payment methods like QR codes or digital wallets."
$wsFunc.Range("C6").Value = "The TVM must be able to accept payment via ATM card or mobile payment"

# 8: "credit/debit card" -> "ATM card" on the Domain Requirements sheet
$wsDomain.Range("E6").Value = "The TVm must support mode of payment such as ATM card, QR code payment and digital wallets"
